$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 9).Value = 5603
$ws.Cells.Item(3, 9).Value = 5782
$ws.Cells.Item(5, 3).Value = 363
$ws.Cells.Item(5, 7).Value = 3708
$ws.Cells.Item(7, 9).Value = 6710
$ws.Cells.Item(10, 9).Value = 7493
$ws.Cells.Item(11, 9).Value = 7882
$ws.Cells.Item(15, 9).Value = 9438
$ws.Cells.Item(18, 7).Value = 5284
$ws.Cells.Item(19, 9).Value = 10770
$ws.Cells.Item(20, 7).Value = 5234
$ws.Cells.Item(22, 7).Value = 5633
$ws.Cells.Item(23, 9).Value = 12155
$ws.Cells.Item(24, 9).Value = 12862
$ws.Cells.Item(26, 9).Value = 13994
$ws.Cells.Item(27, 9).Value = 14543
$ws.Cells.Item(28, 7).Value = 7256
$ws.Cells.Item(28, 9).Value = 14682
$ws.Cells.Item(29, 7).Value = 7522
$ws.Cells.Item(30, 7).Value = 7745
$ws.Cells.Item(30, 9).Value = 15556
$ws.Cells.Item(31, 7).Value = 8034
$ws.Cells.Item(33, 7).Value = 8416
$ws.Cells.Item(33, 9).Value = 17076
$ws.Cells.Item(34, 7).Value = 8503
$ws.Cells.Item(35, 9).Value = 17319
$ws.Cells.Item(36, 7).Value = 8902
$ws.Cells.Item(37, 7).Value = 9142
$ws.Cells.Item(38, 7).Value = 9507
$ws.Cells.Item(39, 7).Value = 9899
$ws.Cells.Item(40, 7).Value = 10177
$ws.Cells.Item(41, 7).Value = 10501
$ws.Cells.Item(41, 9).Value = 20839
$ws.Cells.Item(42, 7).Value = 10969
$ws.Cells.Item(43, 7).Value = 11345
$ws.Cells.Item(45, 7).Value = 12140
$ws.Cells.Item(46, 7).Value = 12520
$ws.Cells.Item(47, 7).Value = 13028
$ws.Cells.Item(48, 7).Value = 13421
$ws.Cells.Item(49, 7).Value = 14072
$ws.Cells.Item(49, 9).Value = 26445
$ws.Cells.Item(50, 7).Value = 14353
$ws.Cells.Item(51, 7).Value = 14468
$ws.Cells.Item(52, 7).Value = 14607
$ws.Cells.Item(53, 7).Value = 14488
$ws.Cells.Item(54, 7).Value = 14496
$ws.Cells.Item(55, 7).Value = 15208
$ws.Cells.Item(55, 9).Value = 28491
$ws.Cells.Item(56, 7).Value = 15620
$ws.Cells.Item(57, 7).Value = 15760
$ws.Cells.Item(58, 7).Value = 15731
$ws.Cells.Item(59, 7).Value = 15962
$ws.Cells.Item(60, 7).Value = 15999
$ws.Cells.Item(60, 9).Value = 30589
$ws.Cells.Item(61, 7).Value = 16450
$ws.Cells.Item(62, 3).Value = 787
$ws.Cells.Item(62, 7).Value = 16592
$ws.Cells.Item(62, 9).Value = 31780
$ws.Cells.Item(63, 7).Value = 16837
$ws.Cells.Item(64, 7).Value = 16668
$ws.Cells.Item(64, 9).Value = 33176
$ws.Cells.Item(66, 7).Value = 17172
$ws.Cells.Item(66, 9).Value = 34268
$ws.Cells.Item(67, 7).Value = 17362
$ws.Cells.Item(69, 7).Value = 17870
$ws.Cells.Item(69, 9).Value = 34902
$ws.Cells.Item(70, 5).Value = 3930
$ws.Cells.Item(70, 7).Value = 17944
$ws.Cells.Item(71, 7).Value = 17248
$ws.Cells.Item(71, 9).Value = 35557
$ws.Cells.Item(72, 7).Value = 17102
$ws.Cells.Item(73, 7).Value = 17230
$ws.Cells.Item(75, 7).Value = 18119
$ws.Cells.Item(76, 7).Value = 19459
$ws.Cells.Item(77, 3).Value = 1061
$ws.Cells.Item(79, 7).Value = 21744
$ws.Cells.Item(81, 5).Value = 5332
$ws.Cells.Item(81, 7).Value = 23818
$ws.Cells.Item(82, 7).Value = 24247
$ws.Cells.Item(82, 9).Value = 42592
$ws.Cells.Item(83, 7).Value = 25079
$ws.Cells.Item(85, 7).Value = 27016
$ws.Cells.Item(86, 7).Value = 28070
$ws.Cells.Item(87, 9).Value = 47695
$ws.Cells.Item(88, 7).Value = 29181
$ws.Cells.Item(89, 7).Value = 30552
$ws.Cells.Item(90, 7).Value = 31484
$ws.Cells.Item(91, 5).Value = 6533
$ws.Cells.Item(91, 7).Value = 31592
$ws.Cells.Item(92, 7).Value = 31495
$ws.Cells.Item(93, 5).Value = 6420
$ws.Cells.Item(93, 7).Value = 33350
$ws.Cells.Item(93, 9).Value = 53867
$ws.Cells.Item(95, 7).Value = 33269
$ws.Cells.Item(96, 7).Value = 33689
$ws.Cells.Item(96, 9).Value = 55299
$ws.Cells.Item(97, 7).Value = 33461
$ws.Cells.Item(98, 7).Value = 34686
$ws.Cells.Item(99, 7).Value = 35044
$ws.Cells.Item(100, 7).Value = 35082
$ws.Cells.Item(101, 7).Value = 35886
$ws.Cells.Item(102, 7).Value = 36495
$ws.Cells.Item(103, 7).Value = 38285
$ws.Cells.Item(104, 7).Value = 40337
$ws.Cells.Item(105, 7).Value = 41805
$ws.Cells.Item(105, 9).Value = 69619
$ws.Cells.Item(106, 7).Value = 42526
$ws.Cells.Item(107, 7).Value = 45140
$ws.Cells.Item(107, 9).Value = 73862
$ws.Cells.Item(108, 7).Value = 45984
$ws.Cells.Item(109, 7).Value = 46335
$ws.Cells.Item(109, 9).Value = 75164
$ws.Cells.Item(110, 7).Value = 47585
$ws.Cells.Item(110, 9).Value = 77173
$ws.Cells.Item(111, 7).Value = 49325
$ws.Cells.Item(112, 7).Value = 50072
$ws.Cells.Item(113, 7).Value = 49929
$ws.Cells.Item(113, 9).Value = 81366
$ws.Cells.Item(114, 7).Value = 50609
$ws.Cells.Item(115, 7).Value = 49951
$ws.Cells.Item(116, 5).Value = 12528
$ws.Cells.Item(117, 7).Value = 52222
$ws.Cells.Item(117, 9).Value = 86556
$ws.Cells.Item(118, 7).Value = 52293
$ws.Cells.Item(118, 9).Value = 87408
$ws.Cells.Item(119, 7).Value = 52975
$ws.Cells.Item(119, 9).Value = 88183
$ws.Cells.Item(120, 7).Value = 54389
$ws.Cells.Item(120, 9).Value = 90866
$ws.Cells.Item(121, 7).Value = 55529
$ws.Cells.Item(121, 9).Value = 93024
$ws.Cells.Item(122, 3).Value = 2848
$ws.Cells.Item(122, 5).Value = 13795
$ws.Cells.Item(122, 7).Value = 56549
$ws.Cells.Item(122, 9).Value = 93578
$ws.Cells.Item(123, 7).Value = 56643
$ws.Cells.Item(123, 9).Value = 94420
$ws.Cells.Item(124, 9).Value = 96482
$ws.Cells.Item(125, 7).Value = 57810
$ws.Cells.Item(125, 9).Value = 97805
$ws.Cells.Item(126, 7).Value = 57783
$ws.Cells.Item(126, 9).Value = 99643
$ws.Cells.Item(127, 9).Value = 100564
$ws.Cells.Item(128, 3).Value = 2964
$ws.Cells.Item(128, 5).Value = 14869
$ws.Cells.Item(128, 7).Value = 59271
$ws.Cells.Item(128, 9).Value = 100723
$ws.Cells.Item(129, 7).Value = 59274
$ws.Cells.Item(129, 9).Value = 100420
$ws.Cells.Item(130, 3).Value = 3064
$ws.Cells.Item(130, 7).Value = 61267
$ws.Cells.Item(130, 9).Value = 101968
$ws.Cells.Item(131, 3).Value = 3012
$ws.Cells.Item(131, 5).Value = 15989
$ws.Cells.Item(131, 7).Value = 62954
$ws.Cells.Item(131, 9).Value = 103959
$ws.Cells.Item(132, 3).Value = 3010
$ws.Cells.Item(132, 5).Value = 16185
$ws.Cells.Item(132, 7).Value = 63146
$ws.Cells.Item(132, 9).Value = 105331
$ws.Cells.Item(133, 3).Value = 3031
$ws.Cells.Item(133, 7).Value = 64036
$ws.Cells.Item(133, 9).Value = 108835
$ws.Cells.Item(134, 3).Value = 3125
$ws.Cells.Item(134, 7).Value = 64391
$ws.Cells.Item(134, 9).Value = 110481
$ws.Cells.Item(135, 3).Value = 3056
$ws.Cells.Item(135, 5).Value = 16893
$ws.Cells.Item(135, 7).Value = 65125
$ws.Cells.Item(135, 9).Value = 112452
$ws.Cells.Item(136, 3).Value = 3072
$ws.Cells.Item(136, 5).Value = 17767
$ws.Cells.Item(136, 7).Value = 66664
$ws.Cells.Item(136, 9).Value = 114345
$ws.Cells.Item(137, 3).Value = 3320
$ws.Cells.Item(137, 5).Value = 19026
$ws.Cells.Item(137, 7).Value = 68606
$ws.Cells.Item(137, 9).Value = 116852
$ws.Cells.Item(138, 3).Value = 3509
$ws.Cells.Item(138, 7).Value = 69453
$ws.Cells.Item(138, 9).Value = 118328
$ws.Cells.Item(139, 3).Value = 3850
$ws.Cells.Item(139, 5).Value = 23007
$ws.Cells.Item(139, 7).Value = 70117
$ws.Cells.Item(139, 9).Value = 119652
$ws.Cells.Item(140, 3).Value = 4837
$ws.Cells.Item(140, 5).Value = 27938
$ws.Cells.Item(140, 7).Value = 70662
$ws.Cells.Item(140, 9).Value = 121118
$ws.Cells.Item(141, 3).Value = 5098
$ws.Cells.Item(141, 5).Value = 28387
$ws.Cells.Item(141, 7).Value = 70192
$ws.Cells.Item(141, 9).Value = 118726
$ws.Cells.Item(142, 3).Value = 5754
$ws.Cells.Item(142, 5).Value = 30531
$ws.Cells.Item(142, 7).Value = 71116
$ws.Cells.Item(142, 9).Value = 119689
$ws.Cells.Item(143, 3).Value = 6252
$ws.Cells.Item(143, 5).Value = 34057
$ws.Cells.Item(143, 7).Value = 72723
$ws.Cells.Item(143, 8).Value = 122459
$ws.Cells.Item(143, 9).Value = 122031
